$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "27.552.43"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.835.44"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.35%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "313.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4240"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.31%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3660"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.07221"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.8646"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "20.73"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.865.38"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.01%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.378"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.515"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.04%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.06954"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "79.84"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.000008980"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "15.42"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "27.849.10"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "5.029"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.65%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.139.83"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.65%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "1.967"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "154.01"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "18.36"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.03%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "5.236"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "114.75"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.94%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "1.818"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.54%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.08875"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.7718"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.90%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "4.550"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.956"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.59%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.150"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.33%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "1.100"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.05361"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.01942"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.818"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.5111"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.56%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "6.838"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.1649"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +1.75%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "10.45"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.39%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.06536"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "105.91"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.4692"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("E49").Value = "  +0.41%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.624"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.796"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.63%  "

